$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D38").Value = "17/6/2025"
$ws.Range("E38").Value = 383
$ws.Range("F38").Value = 543
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 1012
$ws.Range("J38").Value = "N/A"

$ws.Range("D39").Select()
